$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52
$ws.Range("C52").Value = "Giant Hawkfish"
$ws.Range("D52").Value = "Cirrhitus rivulatus"
$ws.Range("G52").Value = 4
$ws.Range("I52").Value = 4
$ws.Range("L52").Value = 4
$ws.Range("M52").Value = 60
$ws.Range("N52").Value = "WhenPlayed"
$ws.Range("O52").Value = "[FishHatch][SchoolFeederMove]"
$ws.Range("S52").Value = 1
$ws.Range("V52").Value = "Named for its hawklike hunting technique, it perches on the high point of a coral reef and dives down upon prey."

# Row 53
$ws.Range("C53").Value = "Giant Manta Ray"
$ws.Range("D53").Value = "Mohula birostris"
$ws.Range("E53").Value = 2
$ws.Range("G53").Value = 1
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = 1
$ws.Range("L53").Value = 5
$ws.Range("M53").Value = 900
$ws.Range("N53").Value = "WhenPlayed"
$ws.Range("O53").Value = "[YoungFish][YoungFish][YoungFish]"
$ws.Range("V53").Value = "The largest ray in the world, it has the biggest brain of any fish by volume—ten times larger than a whale shark’s."

# Row 54
$ws.Range("C54").Value = "Giant Oarfish"
$ws.Range("D54").Value = "Regalecus glesne"
$ws.Range("E54").Value = 1
$ws.Range("G54").Value = 2
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = 1
$ws.Range("L54").Value = 9
$ws.Range("M54").Value = 1100
$ws.Range("N54").Value = "WhenPlayed"
$ws.Range("O54").Value = "(all players) [FishEgg][ArrowDown][FlipperBlue] on each [AllPlayers]"
$ws.Range("S54").Value = 2
$ws.Range("V54").Value = "Rarely seen, the oarfish is the world’s longest bony fish. It’s believed to be the cause of sea serpent sightings."

# Row 55
$ws.Range("C55").Value = "Giant Trevally"
$ws.Range("D55").Value = "Caranx ignobilis"
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("I55").Value = 1
$ws.Range("L55").Value = 5
$ws.Range("M55").Value = 170
$ws.Range("N55").Value = "GameEnd"
$ws.Range("O55").Value = "[FishFromHand][ArrowDown][FlipperPurple]"
$ws.Range("V55").Value = "The giant trevally is known to prey upon fledgling seabirds that fall into the water while learning to fly."

# Row 56
$ws.Range("C56").Value = "Gray Triggerfish"
$ws.Range("D56").Value = "Balistes capriscus"
$ws.Range("E56").Value = 1
$ws.Range("F56").Value = 2
$ws.Range("I56").Value = 1
$ws.Range("L56").Value = 5
$ws.Range("M56").Value = 60
$ws.Range("N56").Value = "IfActivated"
$ws.Range("O56").Value = "[Discard]"
$ws.Range("U56").Value = "blue"
$ws.Range("V56").Value = "Its first, tall dorsal spine remains erect until the smaller second spine is deflexed, triggering the first."

# Row 57
$ws.Range("C57").Value = "Great Northern Tilefish"
$ws.Range("D57").Value = "Lopholatilus chamaeleonticeps"
$ws.Range("E57").Value = 2
$ws.Range("F57").Value = 1
$ws.Range("I57").Value = 1
$ws.Range("J57").Value = 1
$ws.Range("L57").Value = 4
$ws.Range("M57").Value = 125
$ws.Range("N57").Value = "IfActivated"
$ws.Range("O57").Value = "[FishHatch][FishHatch]"
$ws.Range("Q57").Value = 1
$ws.Range("U57").Value = "green"
$ws.Range("V57").Value = "This colorful fish is known as the clown of the sea. It burrows into the sediment at the bottom of the ocean."

# Row 58
$ws.Range("C58").Value = "Great White Shark"
$ws.Range("D58").Value = "Carcharodon carcharias"
$ws.Range("G58").Value = 2
$ws.Range("H58").Value = 1
$ws.Range("I58").Value = 1
$ws.Range("J58").Value = 1
$ws.Range("K58").Value = 1
$ws.Range("L58").Value = 10
$ws.Range("M58").Value = 600
$ws.Range("N58").Value = "WhenPlayed"
$ws.Range("O58").Value = "(all players) [FishEgg][ArrowDown][Predator] on each [AllPlayers]"
$ws.Range("S58").Value = 2
$ws.Range("V58").Value = "Known by scientists as simply the “white shark,” this famous predator is, itself, occasionally preyed upon by orca whales."

# Row 59
$ws.Range("C59").Value = "Greenland Shark"
$ws.Range("D59").Value = "Somniosus microcephalus"
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 1
$ws.Range("J59").Value = 1
$ws.Range("K59").Value = 1
$ws.Range("L59").Value = 8
$ws.Range("M59").Value = 550
$ws.Range("N59").Value = "WhenPlayed"
$ws.Range("O59").Value = "(all players) [FishEgg][ArrowDown][FishLengthLarge] on each [AllPlayers]"
$ws.Range("S59").Value = 2
$ws.Range("V59").Value = "This shark currently holds the record as the longest-lived vertebrate, with an estimated lifespan of over 250 years."

# Row 60
$ws.Range("C60").Value = "Grideye Fish"
$ws.Range("D60").Value = "Ipnops agassizii"
$ws.Range("F60").Value = 2
$ws.Range("K60").Value = 1
$ws.Range("L60").Value = 3
$ws.Range("M60").Value = 16
$ws.Range("N60").Value = "WhenPlayed"
$ws.Range("O60").Value = "[FishEgg][ArrowDown][PlayFishBottomRow] on each"
$ws.Range("V60").Value = "Its eyes are flat, cornea-like, light-sensitive organs without lenses that cover most of the upper surface of the head."

# Row 61
$ws.Range("C61").Value = "Haddock"
$ws.Range("D61").Value = "Melanogrammus aeglefinus"
$ws.Range("E61").Value = 1
$ws.Range("F61").Value = 2
$ws.Range("I61").Value = 1
$ws.Range("J61").Value = 1
$ws.Range("L61").Value = 7
$ws.Range("M61").Value = 110
$ws.Range("N61").Value = "IfActivated"
$ws.Range("O61").Value = "(all players) [DrawCard][AllPlayers]"
$ws.Range("V61").Value = "It lives on gravelly, sandy, and pebbly seafloors, where it feeds on small fish, worms, mollusks, and eggs."

# Row 62
$ws.Range("C62").Value = "Hogfish"
$ws.Range("D62").Value = "Lachnolaimus maximus"
$ws.Range("E62").Value = 2
$ws.Range("G62").Value = 1
$ws.Range("I62").Value = 1
$ws.Range("L62").Value = 6
$ws.Range("M62").Value = 90
$ws.Range("N62").Value = "WhenPlayed"
$ws.Range("O62").Value = "[Discard][Discard][SchoolFeederMove]"
$ws.Range("V62").Value = "Its common name comes from its long pig-like snout and how it roots around the seafloor searching for food."

# Row 63
$ws.Range("C63").Value = "Honeycomb Scaly Dragonfish"
$ws.Range("D63").Value = "Stomias affinis"
$ws.Range("E63").Value = 1
$ws.Range("F63").Value = 2
$ws.Range("I63").Value = 1
$ws.Range("J63").Value = 1
$ws.Range("K63").Value = 1
$ws.Range("L63").Value = 4
$ws.Range("M63").Value = 22
$ws.Range("N63").Value = "GameEnd"
$ws.Range("O63").Value = "[YoungFish][SchoolFeederMove]"
$ws.Range("P63").Value = 1
$ws.Range("S63").Value = 1
$ws.Range("V63").Value = "This fish can be found in deep waters during the day and in shallower areas at night, where it searches for food."

# Row 64
$ws.Range("C64").Value = "Humpback Anglerfish"
$ws.Range("D64").Value = "Melanocetus johnsonii"
$ws.Range("G64").Value = 2
$ws.Range("J64").Value = 1
$ws.Range("K64").Value = 1
$ws.Range("L64").Value = 8
$ws.Range("M64").Value = 18
$ws.Range("N64").Value = "WhenPlayed"
$ws.Range("O64").Value = "(all players) [FishEgg][ArrowDown][PlayFishBottomRow] on each [AllPlayers]"
$ws.Range("P64").Value = 1
$ws.Range("S64").Value = 2
$ws.Range("V64").Value = "The humpback anglerfish is more commonly found at shallower depths than other species in its genus."

# Row 65
$ws.Range("C65").Value = "Humphead Wrasse"
$ws.Range("D65").Value = "Cheilinus undulatus"
$ws.Range("F65").Value = 2
$ws.Range("I65").Value = 1
$ws.Range("L65").Value = 4
$ws.Range("M65").Value = 225
$ws.Range("N65").Value = "WhenPlayed"
$ws.Range("O65").Value = "[FishEgg][ArrowDown][FlipperGreen] on each"
$ws.Range("V65").Value = "It excavates by ejecting water to displace sand. Then it noses around for food."

# Row 66
$ws.Range("C66").Value = "Indo-Pacific Sailfish"
$ws.Range("D66").Value = "Istiophorus platypterus"
$ws.Range("E66").Value = 2
$ws.Range("G66").Value = 1
$ws.Range("I66").Value = 1
$ws.Range("L66").Value = 7
$ws.Range("M66").Value = 350
$ws.Range("N66").Value = "WhenPlayed"
$ws.Range("O66").Value = "[DrawCard][DrawCard][FishHatch]"
$ws.Range("S66").Value = 1
$ws.Range("V66").Value = "The world’s fastest swimming fish, its intricate dorsal fin, shaped like a sail, helps it to reach up to 110 km per hour."

# Row 67
$ws.Range("C67").Value = "Japanese Anchovy"
$ws.Range("D67").Value = "Endraulis japonicus"
$ws.Range("E67").Value = 1
$ws.Range("F67").Value = 1
$ws.Range("I67").Value = 1
$ws.Range("J67").Value = 1
$ws.Range("L67").Value = 3
$ws.Range("M67").Value = 18
$ws.Range("N67").Value = "IfActivated"
$ws.Range("O67").Value = "[FishEgg]"
$ws.Range("U67").Value = "purple"
$ws.Range("V67").Value = "Large schools tend to swim close to the surface during the full moon, their silver scales shimmering in the light."

# Row 68
$ws.Range("C68").Value = "Kaluga"
$ws.Range("D68").Value = "Huso dauricus"
$ws.Range("E68").Value = 2
$ws.Range("G68").Value = 1
$ws.Range("I68").Value = 1
$ws.Range("L68").Value = 6
$ws.Range("M68").Value = 560
$ws.Range("N68").Value = "WhenPlayed"
$ws.Range("O68").Value = "[FishFromHand][ArrowDown][Estuary]"
$ws.Range("V68").Value = "Endangered throughout its range from overfishing, this sturgeon spends much of its life in estuaries and rivers."

# Update the view state: scroll so row 33 is at the top, and select C69
$ws.Range("C69").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
